# Update cryptocurrency price/volume data per the Wed Jul 5 19:11:50 UTC 2023 GitHub Actions run
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '30.513.44'
$ws.Range("E2").Value = '  -1.25%  '

# Row 3
$ws.Range("D3").Value = '1.911.45'
$ws.Range("E3").Value = '  -1.94%  '

# Row 4
$ws.Range("E4").Value = '  -0.09%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '239.62'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.94%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4783'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -1.81%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2843'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -2.79%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06686'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -2.38%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '18.65'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -3.94%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '100.93'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -4.83%  '

# Row 12
$ws.Range("B12").Value = 'TRON'
$ws.Range("C12").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07688'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -0.67%  '

# Row 13
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.908.05'
$ws.Range("E13").Value = '  -2.09%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.221'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -1.57%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6686'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -3.57%  '

# Row 16
$ws.Range("D16").Value = '30.506.58'
$ws.Range("E16").Value = '  -1.28%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '255.65'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -8.00%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.001'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +0.04%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007465'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -2.99%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.66'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -3.68%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.377'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -0.97%  '

# Row 22
$ws.Range("E22").Value = '  -0.07%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.287'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -2.68%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.317'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -3.79%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '166.93'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -0.28%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '19.05'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -2.82%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.055'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -4.74%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '4.714'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +3.70%  '

# Row 29
$ws.Range("E29").Value = '  -2.58%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.381'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -0.91%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.513'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -2.60%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.254'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -2.40%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04716'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -2.54%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7293'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -2.21%  '

# Row 35
$ws.Range("E35").Value = '  -4.14%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.9995'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -0.02%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.703'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -1.04%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01916'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -3.45%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.609'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -2.61%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '74.80'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -2.27%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.221'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -3.85%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.968'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -5.55%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.8619'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -3.88%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '105.21'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -2.50%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.001'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +0.15%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4239'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -3.68%  '

# Row 47
$ws.Range("E47").Value = '  -3.68%  '

# Row 48
$ws.Range("B48").Value = 'Algorand'
$ws.Range("C48").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.1198'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -3.52%  '

# Row 49
$ws.Range("B49").Value = 'Elrond'
$ws.Range("C49").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '34.71'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -2.50%  '

# Row 50
$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.773'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -4.32%  '

# Row 51
$ws.Range("B51").Value = 'Maker'
$ws.Range("C51").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '890.54'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -10.24%  '
